# Applies the diff: rows 8<->9 swap, rows 15<->16 swap, and rows 19->20->21->19
# rotate (row19 takes old row20 data, row20 takes old row21 data, row21 takes
# old row19 data). Only the fields that actually differ are touched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8 <-> Row 9 ------------------------------------------------------
$ws.Range("A8").Value = 130930217
$ws.Range("B8").Value = 79499
$ws.Range("D8").Value = "LC"
$ws.Range("E8").Value = 6459
$ws.Range("F8").Value = "Barkkornlav"
$ws.Range("G8").Value = "Lopadium disciforme"
$ws.Range("H8").Value = "(Flot.) Kullh."
$ws.Range("Q8").Value = 448392
$ws.Range("R8").Value = 7037298
$ws.Range("AH8").Value = ""
$ws.Range("AJ8").Value = "rönn"
$ws.Range("AK8").Value = "Sorbus aucuparia"
$ws.Range("AO8").Value = "Sorbus aucuparia"

$ws.Range("A9").Value = 130930226
$ws.Range("B9").Value = 91828
$ws.Range("D9").Value = "NT"
$ws.Range("E9").Value = 5432
$ws.Range("F9").Value = "Granticka"
$ws.Range("G9").Value = "Porodaedalea chrysoloma s.lat."
$ws.Range("H9").Value = ""
$ws.Range("Q9").Value = 448362
$ws.Range("R9").Value = 7037345
$ws.Range("AH9").Value = "Granskog"
$ws.Range("AJ9").Value = "gran"
$ws.Range("AK9").Value = "Picea abies"
$ws.Range("AO9").Value = "Picea abies"

# --- Row 15 <-> Row 16 -----------------------------------------------------
$ws.Range("A15").Value = 130930220
$ws.Range("B15").Value = 79714
$ws.Range("D15").Value = "NT"
$ws.Range("E15").Value = 1797
$ws.Range("F15").Value = "Mjölig dropplav"
$ws.Range("G15").Value = "Cliostomum leprosum"
$ws.Range("H15").Value = "(Räsänen) Holien & Tønsberg"
$ws.Range("Q15").Value = 448353
$ws.Range("R15").Value = 7037267

$ws.Range("A16").Value = 130930219
$ws.Range("B16").Value = 92530
$ws.Range("D16").Value = "LC"
$ws.Range("E16").Value = 3298
$ws.Range("F16").Value = "Trådticka"
$ws.Range("G16").Value = "Climacocystis borealis"
$ws.Range("H16").Value = "(Fr.) Kotl. & Pouzar"
$ws.Range("Q16").Value = 448355
$ws.Range("R16").Value = 7037273

# --- Row 19 -> Row 20 -> Row 21 -> Row 19 (3-way rotation) -----------------
$ws.Range("A19").Value = 130930231
$ws.Range("B19").Value = 83223
$ws.Range("D19").Value = "NT"
$ws.Range("E19").Value = 6440
$ws.Range("F19").Value = "Vitgrynig nållav"
$ws.Range("G19").Value = "Chaenotheca subroscida"
$ws.Range("H19").Value = "(Eitner) Zahlbr."
$ws.Range("Q19").Value = 448412
$ws.Range("R19").Value = 7037419

$ws.Range("A20").Value = 130930223
$ws.Range("B20").Value = 79714
$ws.Range("E20").Value = 1797
$ws.Range("F20").Value = "Mjölig dropplav"
$ws.Range("G20").Value = "Cliostomum leprosum"
$ws.Range("H20").Value = "(Räsänen) Holien & Tønsberg"
$ws.Range("Q20").Value = 448337
$ws.Range("R20").Value = 7037328

$ws.Range("A21").Value = 130930222
$ws.Range("B21").Value = 83221
$ws.Range("D21").Value = "VU"
$ws.Range("E21").Value = 6486
$ws.Range("F21").Value = "Skuggnål"
$ws.Range("G21").Value = "Chaenotheca sphaerocephala"
$ws.Range("H21").Value = "Nádv."
$ws.Range("Q21").Value = 448330
$ws.Range("R21").Value = 7037323
